# Generate Report for Handoff
# Adds a new tracked file (ef4783de-cfb1-4ba8-ade8-09b399ca2752.md) to the
# localization-status workbook: one summary row on "Overview" and one
# detail row on each of the "zh-cn" / "de-de" language sheets.

$wb = $excel.ActiveWorkbook

$newGuid       = "ef4783de-cfb1-4ba8-ade8-09b399ca2752"
$newGuidMd     = $newGuid + ".md"
$newHash       = "69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6"

$mdUrl         = "https://github.com/OpenLocalizationTest/oltest/blob/4daf793447831388b2108a25df80716aaad753b6/e2e/" + $newGuidMd

# ---------------------------------------------------------------------
# Sheet "Overview": append summary row for the new file
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = $newGuidMd
$ws1.Range("A3").Style = "HyperLink"
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdUrl, "", "", $newGuidMd) | Out-Null

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-24-11 16:24:31"

# ---------------------------------------------------------------------
# Sheet "zh-cn": append detail row for the new file
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$zhXlf    = $newGuid + "." + $newHash + ".zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dabd5264373e6c580524dcffc1c65dc479913ba3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $zhXlf

$ws2.Range("A3").Value = $newGuidMd
$ws2.Range("A3").Style = "HyperLink"
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdUrl, "", "", $newGuidMd) | Out-Null

$ws2.Range("B3").Value = ".md"
$ws2.Range("B3").Style = "HyperLink"
$ws2.Hyperlinks.Add($ws2.Range("B3"), $mdUrl, "", "", ".md") | Out-Null

$ws2.Range("C3").Value = "Ready for handoff"

$ws2.Range("D3").Value = $zhXlf
$ws2.Range("D3").Style = "HyperLink"
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhXlfUrl, "", "", $zhXlf) | Out-Null

$ws2.Range("E3").Value = "2016-03-11 16:24:28"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de": append detail row for the new file
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$deXlf    = $newGuid + "." + $newHash + ".de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e328e07d1d5610e0ec79e93dab15b6fbd640fd09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $deXlf

$ws3.Range("A3").Value = $newGuidMd
$ws3.Range("A3").Style = "HyperLink"
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdUrl, "", "", $newGuidMd) | Out-Null

$ws3.Range("B3").Value = ".md"
$ws3.Range("B3").Style = "HyperLink"
$ws3.Hyperlinks.Add($ws3.Range("B3"), $mdUrl, "", "", ".md") | Out-Null

$ws3.Range("C3").Value = "Ready for handoff"

$ws3.Range("D3").Value = $deXlf
$ws3.Range("D3").Style = "HyperLink"
$ws3.Hyperlinks.Add($ws3.Range("D3"), $deXlfUrl, "", "", $deXlf) | Out-Null

$ws3.Range("E3").Value = "2016-03-11 16:24:31"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

Write-Host "Handoff report rows added for $newGuidMd"
